# Add season-record columns (Wins / Losses / Ties) to the right of the
# existing team-stats table.
#
# The sheet currently uses columns A:AC (header row 1, data rows 2-56).
# We append three new columns: AD = Wins, AE = Losses, AF = Ties, with the
# same season-record values (60 wins, 102 losses, 0 ties) repeated for
# every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - mirror the formatting of the existing header cells
# (e.g. AC1, bold + bordered + centered) by copying its style onto the new
# header cells before setting their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$wins = 60
$losses = 102
$ties = 0

for ($row = 2; $row -le 56; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins    # column AD = 30
    $ws.Cells.Item($row, 31).Value = $losses  # column AE = 31
    $ws.Cells.Item($row, 32).Value = $ties    # column AF = 32
}
